$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep text formatting so numeric-looking strings
# (e.g. "1.001") are not auto-converted to numbers by Excel value parsing.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.852.50'
$ws.Range("E2").Value = '  -1.72%  '

$ws.Range("D3").Value = '1.806.09'
$ws.Range("E3").Value = '  -1.21%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.48%  '

$ws.Range("D5").Value = '309.90'
$ws.Range("E5").Value = '  -1.45%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("D7").Value = '0.4424'
$ws.Range("E7").Value = '  +4.16%  '

$ws.Range("D8").Value = '0.3663'
$ws.Range("E8").Value = '  -1.27%  '

$ws.Range("D9").Value = '0.07308'
$ws.Range("E9").Value = '  +0.64%  '

$ws.Range("D10").Value = '0.8541'
$ws.Range("E10").Value = '  -1.44%  '

$ws.Range("D11").Value = '20.64'
$ws.Range("E11").Value = '  -2.27%  '

$ws.Range("D12").Value = '1.806.88'
$ws.Range("E12").Value = '  -1.11%  '

$ws.Range("D13").Value = '6.599'
$ws.Range("E13").Value = '  -2.05%  '

$ws.Range("D14").Value = '0.07099'
$ws.Range("E14").Value = '  +0.04%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '92.17'
$ws.Range("E15").Value = '  +2.82%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '5.306'
$ws.Range("E16").Value = '  -0.29%  '

$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.51%  '

$ws.Range("D18").Value = '0.000008700'
$ws.Range("E18").Value = '  -1.98%  '

$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("D20").Value = '14.86'
$ws.Range("E20").Value = '  -1.60%  '

$ws.Range("D21").Value = '26.869.84'
$ws.Range("E21").Value = '  -2.09%  '

$ws.Range("D22").Value = '5.155'
$ws.Range("E22").Value = '  +0.46%  '

$ws.Range("D23").Value = '10.84'
$ws.Range("E23").Value = '  -0.67%  '

$ws.Range("E24").Value = '  -0.52%  '

$ws.Range("D25").Value = '151.47'
$ws.Range("E25").Value = '  -1.02%  '

$ws.Range("D26").Value = '2.204'
$ws.Range("E26").Value = '  +1.35%  '

$ws.Range("D27").Value = '18.52'
$ws.Range("E27").Value = '  +0.55%  '

$ws.Range("D28").Value = '5.197'
$ws.Range("E28").Value = '  -0.85%  '

$ws.Range("D29").Value = '116.74'
$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("D30").Value = '0.08793'
$ws.Range("E30").Value = '  -0.98%  '

$ws.Range("D31").Value = '1.176'
$ws.Range("E31").Value = '  -1.97%  '

$ws.Range("D32").Value = '0.7493'
$ws.Range("E32").Value = '  -1.22%  '

$ws.Range("D33").Value = '2.930'
$ws.Range("E33").Value = '  +3.72%  '

$ws.Range("D34").Value = '4.447'
$ws.Range("E34").Value = '  -0.39%  '

$ws.Range("D35").Value = '1.000'
$ws.Range("E35").Value = '  -0.51%  '

$ws.Range("E36").Value = '  -2.45%  '

$ws.Range("D37").Value = '0.01965'
$ws.Range("E37").Value = '  -0.84%  '

$ws.Range("D38").Value = '0.05190'
$ws.Range("E38").Value = '  -1.54%  '

$ws.Range("D39").Value = '0.5369'
$ws.Range("E39").Value = '  +5.95%  '

$ws.Range("D40").Value = '2.864'
$ws.Range("E40").Value = '  -0.24%  '

$ws.Range("D41").Value = '7.037'
$ws.Range("E41").Value = '  -4.74%  '

$ws.Range("D42").Value = '0.1689'
$ws.Range("E42").Value = '  -0.73%  '

$ws.Range("D43").Value = '0.5219'
$ws.Range("E43").Value = '  +9.50%  '

$ws.Range("D44").Value = '8.432'
$ws.Range("E44").Value = '  -3.15%  '

$ws.Range("D45").Value = '10.57'
$ws.Range("E45").Value = '  -0.49%  '

$ws.Range("D46").Value = '1.972'
$ws.Range("E46").Value = '  +6.04%  '

$ws.Range("D47").Value = '105.45'
$ws.Range("E47").Value = '  -2.18%  '

$ws.Range("D48").Value = '0.9998'
$ws.Range("E48").Value = '  -0.51%  '

$ws.Range("D49").Value = '1.663'
$ws.Range("E49").Value = '  -0.62%  '

$ws.Range("D50").Value = '0.06331'
$ws.Range("E50").Value = '  -0.89%  '

$ws.Range("D51").Value = '0.9179'
$ws.Range("E51").Value = '  +0.00%  '
